# AFDP-3021: Configuring assignee per queue
# Adds a new "Set Assignee" action column (G) to the "Save Case File Rules"
# rule table on Sheet1, mirroring the existing "Set Owning Group" column (F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formats from column F (the existing last action column) onto the new
# column G for the whole rule-table range so borders/fills/fonts match.
$ws.Range("F19:F34").Copy()
$ws.Range("G19:G34").PasteSpecial(-4122)

# Header cells for the new column.
$ws.Range("G20").Value = "ACTION"
$ws.Range("G22").Value = 'addOrUpdateParticipant($caseFile, "$1", "$2");'
$ws.Range("G23").Value = "Set Assignee"

# Per-queue assignee values.
$ws.Range("G30").Value = "assignee, sally-acm"
$ws.Range("G31").Value = "assignee, ian-acm"
$ws.Range("G32").Value = "assignee, samuel-acm"
$ws.Range("G33").Value = "assignee, ann-acm"
$ws.Range("G34").Value = "assignee, sally-acm"

# Match the column width used for the other wide action columns.
$ws.Columns.Item(7).ColumnWidth = 43.140625

# Reflect the new column in the view/selection state.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("G23").Select()
